$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 864.5
$ws.Range("I29").Value = 194.66667
$ws.Range("J29").Value = 2874
$ws.Range("K29").Value = 584.00001
$ws.Range("L29").Value = 8622
$ws.Range("M29").Value = -303.00001
$ws.Range("N29").Value = -9184
$ws.Range("H38").Value = 45.363636
$ws.Range("I38").Value = 45.363636
$ws.Range("K38").Value = 136.090908
$ws.Range("M38").Value = 235.909092
$ws.Range("H40").Value = 2296.2
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 2420.25
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 2420.25
$ws.Range("M40").Value = -1625
$ws.Range("N40").Value = -2770.25
$ws.Range("H58").Value = 1456.0769
$ws.Range("J58").Value = 1382.3334
$ws.Range("L58").Value = 4147.0002
$ws.Range("N58").Value = -4447.0002
$ws.Range("H87").Value = 39999
$ws.Range("J87").Value = 39999
$ws.Range("L87").Value = 39999
$ws.Range("N87").Value = -42495
$ws.Range("H90").Value = 39999
$ws.Range("J90").Value = 39999
$ws.Range("L90").Value = 119997
$ws.Range("N90").Value = -132477
$ws.Range("H103").Value = 1000
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H107").Value = 1232.5
$ws.Range("I107").Value = 1232.5
$ws.Range("K107").Value = 1232.5
$ws.Range("M107").Value = 687.5
$ws.Range("H132").Value = 1261.8695
$ws.Range("I132").Value = 914.1177
$ws.Range("J132").Value = 2247.1667
$ws.Range("K132").Value = 2742.3531
$ws.Range("L132").Value = 6741.500100000001
$ws.Range("M132").Value = -212.3531000000003
$ws.Range("N132").Value = -11801.5001
$ws.Range("H138").Value = 1613.8334
$ws.Range("I138").Value = 455.33334
$ws.Range("K138").Value = 1366.00002
$ws.Range("M138").Value = 3773.99998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H61").Value = 2039.5
$ws.Range("I61").Value = 1988.4445
$ws.Range("K61").Value = 1988.4445
$ws.Range("M61").Value = -1776.4445
$ws.Range("H63").Value = 2869.2856
$ws.Range("I63").Value = 1795.3334
$ws.Range("J63").Value = 3674.75
$ws.Range("K63").Value = 1795.3334
$ws.Range("L63").Value = 3674.75
$ws.Range("M63").Value = -1109.3334
$ws.Range("N63").Value = -5046.75
$ws.Range("H66").Value = 2869.2856
$ws.Range("I66").Value = 1795.3334
$ws.Range("J66").Value = 3674.75
$ws.Range("K66").Value = 8976.666999999999
$ws.Range("L66").Value = 18373.75
$ws.Range("M66").Value = -5544.666999999999
$ws.Range("N66").Value = -25237.75
$ws.Range("H132").Value = 2379.6667
$ws.Range("I132").Value = 2427.125
$ws.Range("K132").Value = 7281.375
$ws.Range("M132").Value = -4751.375
$ws.Range("H136").Value = 2039.5
$ws.Range("I136").Value = 1988.4445
$ws.Range("K136").Value = 5965.333500000001
$ws.Range("M136").Value = -3415.333500000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H82").Value = 31142.857
$ws.Range("I82").Value = 9000
$ws.Range("K82").Value = 9000
$ws.Range("M82").Value = -8617
$ws.Range("H85").Value = 31142.857
$ws.Range("I85").Value = 9000
$ws.Range("K85").Value = 9000
$ws.Range("M85").Value = -7674
$ws.Range("H107").Value = 1260.2222
$ws.Range("J107").Value = 1900
$ws.Range("L107").Value = 1900
$ws.Range("N107").Value = -5740
$ws.Range("H134").Value = 9200.091
$ws.Range("I134").Value = 9400.125
$ws.Range("K134").Value = 28200.375
$ws.Range("M134").Value = -25665.375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.27273
$ws.Range("I7").Value = 88
$ws.Range("J7").Value = 12.5
$ws.Range("K7").Value = 88
$ws.Range("L7").Value = 12.5
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = -238.5
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H132").Value = 1413.5454
$ws.Range("I132").Value = 1510.1111
$ws.Range("J132").Value = 979
$ws.Range("K132").Value = 4530.3333
$ws.Range("L132").Value = 2937
$ws.Range("M132").Value = -2000.3333
$ws.Range("N132").Value = -7997

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2277.5
$ws.Range("I4").Value = 2188.889
$ws.Range("K4").Value = 6566.667
$ws.Range("M4").Value = -6454.667
$ws.Range("H7").Value = 265.5
$ws.Range("I7").Value = 269.83334
$ws.Range("J7").Value = 252.5
$ws.Range("K7").Value = 809.5000200000001
$ws.Range("L7").Value = 757.5
$ws.Range("M7").Value = -697.5000200000001
$ws.Range("N7").Value = -981.5
$ws.Range("H10").Value = 18.866667
$ws.Range("I10").Value = 18.866667
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 56.600001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 82.39999900000001
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 669333
$ws.Range("J11").Value = 3999.5
$ws.Range("L11").Value = 11998.5
$ws.Range("N11").Value = -12278.5
$ws.Range("H16").Value = 506.14285
$ws.Range("I16").Value = 47.666668
$ws.Range("J16").Value = 850
$ws.Range("K16").Value = 143.000004
$ws.Range("L16").Value = 2550
$ws.Range("M16").Value = 29.99999600000001
$ws.Range("N16").Value = -2896
$ws.Range("H17").Value = 100
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10423950
$ws.Range("I122").Value = 20843066
$ws.Range("K122").Value = 62529198
$ws.Range("M122").Value = -62526748

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9000
$ws.Range("I61").Value = 9000
$ws.Range("K61").Value = 9000
$ws.Range("M61").Value = -8798
$ws.Range("H100").Value = 367
$ws.Range("I100").Value = 367
$ws.Range("K100").Value = 367
$ws.Range("M100").Value = 174
$ws.Range("H113").Value = 9000
$ws.Range("I113").Value = 9000
$ws.Range("K113").Value = 9000
$ws.Range("M113").Value = -6830

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16393.666
$ws.Range("I41").Value = 16473.2
$ws.Range("K41").Value = 16473.2
$ws.Range("M41").Value = -16083.2
$ws.Range("H113").Value = 9157.538
$ws.Range("I113").Value = 10395.546
$ws.Range("K113").Value = 31186.638
$ws.Range("M113").Value = -29016.638
$ws.Range("H132").Value = 1080.0834
$ws.Range("I132").Value = 1080.0834
$ws.Range("K132").Value = 3240.2502
$ws.Range("M132").Value = -710.2501999999999
